# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Swap the two shared-string entries "Polinesia Francesa" / "Belice" ---
# Row 193 used to be "Polinesia Francesa" and row 194 used to be "Belice";
# after the edit their names (and per-country stats) trade places.
$ws.Range("A193").Value = "Belice"
$ws.Range("A194").Value = "Polinesia Francesa"

# --- Refresh the "last updated" timestamp string ---
$ws.Range("A1").Value = "Datos actualizados a 5 de Agosto de 2020 a las 06:48"

# --- Updated per-country case numbers (B:Casos totales, C:Nuevos casos,
#     D:Casos activos, E:Recuperados, F:Casos criticos, G:Muertes hoy, H:Muertes) ---

# India (row 6)
$ws.Range("B6").Value = 1908254
$ws.Range("C6").Value = 1641
$ws.Range("D6").Value = 1282215
$ws.Range("E6").Value = 586219

# Haiti (row 92)
$ws.Range("B92").Value = 7532
$ws.Range("C92").Value = 21
$ws.Range("E92").Value = 2529
$ws.Range("G92").Value = 5
$ws.Range("H92").Value = 171

# Tailandia (row 115)
$ws.Range("B115").Value = 3328
$ws.Range("C115").Value = 7
$ws.Range("E115").Value = 128

# Birmania (row 169)
$ws.Range("B169").Value = 356
$ws.Range("C169").Value = 1
$ws.Range("E169").Value = 48

# Mongolia (row 172)
$ws.Range("D172").Value = 244
$ws.Range("E172").Value = 49

# Butan (row 190)
$ws.Range("B190").Value = 105
$ws.Range("C190").Value = 2
$ws.Range("D190").Value = 93
$ws.Range("E190").Value = 12

# Row 193 (now "Belice")
$ws.Range("B193").Value = 72
$ws.Range("C193").Value = 15
$ws.Range("D193").Value = 31
$ws.Range("E193").Value = 39
$ws.Range("F193").Value = 0
$ws.Range("G193").Value = 0
$ws.Range("H193").Value = 2

# Row 194 (now "Polinesia Francesa")
$ws.Range("B194").Value = 64
$ws.Range("C194").Value = 0
$ws.Range("D194").Value = 62
$ws.Range("E194").Value = 2
$ws.Range("F194").Value = 0
$ws.Range("G194").Value = 0
$ws.Range("H194").Value = 0
